$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1.134434208527111
$ws.Range("F2").Value = 0.455877884661925
$ws.Range("G2").Value = 0.1574095222632352
$ws.Range("H2").Value = 0.184749521947265
$ws.Range("I2").Value = 0.9351544332584775

$ws.Range("E3").Value = 1.136118914357281
$ws.Range("F3").Value = 0.455882018360749
$ws.Range("G3").Value = 0.1457815837570197
$ws.Range("H3").Value = 0.1798919439769675
$ws.Range("I3").Value = 0.918037859851653

$ws.Range("E4").Value = 1.136682560272896
$ws.Range("F4").Value = 0.4558875315652311
$ws.Range("G4").Value = 0.1415897291679757
$ws.Range("H4").Value = 0.1731594211844923
$ws.Range("I4").Value = 0.9053366706458894

$ws.Range("E5").Value = 1.137462449876291
$ws.Range("F5").Value = 0.4558802205251314
$ws.Range("G5").Value = 0.1295884230354892
$ws.Range("H5").Value = 0.2377048333801086
$ws.Range("I5").Value = 0.8783269748364995

$ws.Range("E6").Value = 1.137215553554008
$ws.Range("F6").Value = 0.4558831007021503
$ws.Range("G6").Value = 0.1340592719927949
$ws.Range("H6").Value = 0.2138600778478064
$ws.Range("I6").Value = 0.8869253444704601

$ws.Range("E7").Value = 1.137706560586925
$ws.Range("F7").Value = 0.4558579720385676
$ws.Range("G7").Value = 0.1343193687919772
$ws.Range("H7").Value = 0.2191309382299363
$ws.Range("I7").Value = 0.8929689733717709

$ws.Range("E8").Value = 1.13728077563392
$ws.Range("F8").Value = 0.4559031713540883
$ws.Range("G8").Value = 0.1408291341884509
$ws.Range("H8").Value = 0.2025227718341068
$ws.Range("I8").Value = 0.9142092367328025

$ws.Range("E9").Value = 1.13826111825791
$ws.Range("F9").Value = 0.4558808322025088
$ws.Range("G9").Value = 0.1353228446794389
$ws.Range("H9").Value = 0.2124479991246679
$ws.Range("I9").Value = 0.8935099677788353

$ws.Range("E10").Value = 1.137000388450454
$ws.Range("F10").Value = 0.4558842820522531
$ws.Range("G10").Value = 0.1430123778886961
$ws.Range("H10").Value = 0.1664568728568378
$ws.Range("I10").Value = 0.9179121544253221

$ws.Range("E11").Value = 1.133573279502847
$ws.Range("F11").Value = 0.4558824606103568
$ws.Range("G11").Value = 0.160463463870248
$ws.Range("H11").Value = 0.1246647451358507
$ws.Range("I11").Value = 0.9265988014464368

$ws.Range("E12").Value = 1.135095788030593
$ws.Range("F12").Value = 0.4558822478561606
$ws.Range("G12").Value = 0.1569188596736608
$ws.Range("H12").Value = 0.1443141942799877
$ws.Range("I12").Value = 0.9362895668690131

$ws.Range("E13").Value = 1.13806327037512
$ws.Range("F13").Value = 0.4558827905558146
$ws.Range("G13").Value = 0.1434022503157772
$ws.Range("H13").Value = 0.1763477312243915
$ws.Range("I13").Value = 0.920188573960699

$ws.Range("E14").Value = 1.137277719151322
$ws.Range("F14").Value = 0.4558810404526684
$ws.Range("G14").Value = 0.1448174610017762
$ws.Range("H14").Value = 0.1842467028154326
$ws.Range("I14").Value = 0.9313249555254248

$ws.Range("E15").Value = 1.137085861124306
$ws.Range("F15").Value = 0.4558829401586367
$ws.Range("G15").Value = 0.1452840388754392
$ws.Range("H15").Value = 0.1961669919200982
$ws.Range("I15").Value = 0.9330080285341689

$ws.Range("E16").Value = 1.137631515809276
$ws.Range("F16").Value = 0.4558855760120344
$ws.Range("G16").Value = 0.1420477598290924
$ws.Range("H16").Value = 0.2182465780101087
$ws.Range("I16").Value = 0.9285701793552507

$ws.Range("E17").Value = 1.139005126867447
$ws.Range("F17").Value = 0.4558819955097421
$ws.Range("G17").Value = 0.1374398636605355
$ws.Range("H17").Value = 0.220881263207821
$ws.Range("I17").Value = 0.9137200334139748

$ws.Range("E18").Value = 1.139478091203545
$ws.Range("F18").Value = 0.4558811198441085
$ws.Range("G18").Value = 0.1367541752163105
$ws.Range("H18").Value = 0.2173944970758712
$ws.Range("I18").Value = 0.9091924839029113

